# Updating results for BM, EX and LIN
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# EX sheet: fill in TZFE/OPFE mean & RMSE results (B3:E10)
# ---------------------------------------------------------------
$exws = $wb.Worksheets.Item("EX")
$exData = @(
  @(0.0289612,    0.400054,    0.0683016,    0.507326),
  @(0.00958103,   0.132017,    0.00725576,   0.00909482),
  @(0.0155158,    0.213793,    0.0117502,    0.0147285),
  @(0.000140764,  0.0118049,   0.00394885,   0.00494974),
  @(0.191621,     2.64034,     0.145115,     0.181896),
  @(0.015435,     0.212663,    0.0116856,    0.0146477),
  @(0.0252643,    0.348143,    0.00717294,   0.00911457),
  @(0.00452605,   0.0648859,   0.00128636,   0.00163585)
)
for ($i = 0; $i -lt $exData.Count; $i++) {
  $r = 3 + $i
  for ($j = 0; $j -lt 4; $j++) {
    $c = 2 + $j
    $exws.Cells.Item($r, $c).Value = $exData[$i][$j]
  }
}
$exws.Range("B3:E10").NumberFormat = "0.00%"
$exws.Range("F22").Select()

# ---------------------------------------------------------------
# LIN sheet: same results, shifted one column to the right (C3:F10)
# ---------------------------------------------------------------
$linws = $wb.Worksheets.Item("LIN")
$linData = $exData
for ($i = 0; $i -lt $linData.Count; $i++) {
  $r = 3 + $i
  for ($j = 0; $j -lt 4; $j++) {
    $c = 3 + $j
    $linws.Cells.Item($r, $c).Value = $linData[$i][$j]
  }
}
$linws.Range("C3:F10").NumberFormat = "0.00%"
$linws.Range("C3:F10").Select()
$linws.Range("C3").Activate()

# ---------------------------------------------------------------
# GSSA sheet: not yet solved - only the saved selection moved
# ---------------------------------------------------------------
$gssaws = $wb.Worksheets.Item("GSSA")
$gssaws.Range("H22").Select()

# ---------------------------------------------------------------
# Table sheet: pull TZFE mean/RMSE from EX & LIN, and refresh the
# Euler-error / timing numbers
# ---------------------------------------------------------------
$tws = $wb.Worksheets.Item("Table")

$tws.Range("B3").Formula = "=EX!B3"
$tws.Range("C3").Formula = "=LIN!C3"
$tws.Range("D3").Value = $null
$tws.Range("E3").Value = $null

$tws.Range("B4").Formula = "=EX!C3"
$tws.Range("C4").Formula = "=LIN!D3"
$tws.Range("D4").Value = $null
$tws.Range("E4").Value = $null

$tws.Range("B5").Value = $null
$tws.Range("C5").Value = $null
$tws.Range("D5").Value = $null
$tws.Range("E5").Value = $null

$tws.Range("B6").Formula = "=EX!D3"
$tws.Range("C6").Formula = "=LIN!E3"
$tws.Range("D6").Value = $null
$tws.Range("E6").Value = $null

$tws.Range("B7").Formula = "=EX!E3"
$tws.Range("C7").Formula = "=LIN!F3"
$tws.Range("D7").Value = $null
$tws.Range("E7").Value = $null

$tws.Range("B3:E7").NumberFormat = "0.00%"

# Euler errors / RMSE row
$tws.Range("B9").Value = 0.000242552
$tws.Range("C9").Value = 0.00025194099999999999
$tws.Range("A9:C9").NumberFormat = "0.00E+00"

# Time row
$tws.Range("B10").Value = $null
$tws.Range("C10").Value = $null
$tws.Range("B10:C10").NumberFormat = "0.0000"

$tws.Range("B11").Value = "n/a"
$tws.Range("C11").Value = 0.0018518824345174754
$tws.Range("B11:C11").NumberFormat = "0.0000"
$tws.Range("B11").HorizontalAlignment = -4108

$tws.Range("B12").Value = 31.537829729678734
$tws.Range("C12").Value = 62.380262520635483
$tws.Range("B12:C12").NumberFormat = "0.0000"

$tws.Range("D21").Select()

$excel.ActiveWorkbook.Worksheets.Item("Table").Activate()
